# Update stock report rows: several rows that share the same product name
# had their Opening/Closing/Qty/Value figures (columns B, C, D, E, F, G)
# rearranged among themselves. Apply the corrected values cell-by-cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B161").Value = 64350
$ws.Range("E161").Value = 70.63
$ws.Range("F161").Value = 101
$ws.Range("G161").Value = 6710.44

$ws.Range("B162").Value = 57756
$ws.Range("E162").Value = 79.37
$ws.Range("F162").Value = -100
$ws.Range("G162").Value = -6644

$ws.Range("B163").Value = 53925
$ws.Range("F163").Value = 1
$ws.Range("G163").Value = 66.44

$ws.Range("B183").Value = 64329
$ws.Range("E183").Value = 128.32
$ws.Range("F183").Value = 6
$ws.Range("G183").Value = 724.14

$ws.Range("B184").Value = 57552
$ws.Range("E184").Value = 136.86
$ws.Range("F184").Value = -5
$ws.Range("G184").Value = -603.45

$ws.Range("B264").Value = 48719
$ws.Range("E264").Value = 353.35
$ws.Range("F264").Value = -81
$ws.Range("G264").Value = -23955.75

$ws.Range("B265").Value = 64979
$ws.Range("E265").Value = 314.41
$ws.Range("F265").Value = 82
$ws.Range("G265").Value = 24251.5

$ws.Range("B313").Value = 57854
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.6799999999999

$ws.Range("B314").Value = 62997
$ws.Range("F314").Value = 72
$ws.Range("G314").Value = 22020.48

$ws.Range("B316").Value = 61610
$ws.Range("D316").Value = 102.71
$ws.Range("E316").Value = 122.71
$ws.Range("F316").Value = -58
$ws.Range("G316").Value = -5957.18

$ws.Range("B318").Value = 57077
$ws.Range("D318").Value = 93.08
$ws.Range("E318").Value = 111.2
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 93.08

$ws.Range("B355").Value = 55356
$ws.Range("E355").Value = 54.04
$ws.Range("F355").Value = -158
$ws.Range("G355").Value = -7527.12

$ws.Range("B356").Value = 63510
$ws.Range("E356").Value = 50.66
$ws.Range("F356").Value = 167
$ws.Range("G356").Value = 7955.88

$ws.Range("B372").Value = 63652
$ws.Range("E372").Value = 55.42
$ws.Range("F372").Value = 250
$ws.Range("G372").Value = 13032.5

$ws.Range("B373").Value = 57885
$ws.Range("E373").Value = 62.28
$ws.Range("F373").Value = 4
$ws.Range("G373").Value = 208.52

$ws.Range("B382").Value = 63560
$ws.Range("E382").Value = 134.87
$ws.Range("F382").Value = 104
$ws.Range("G382").Value = 13193.44

$ws.Range("B383").Value = 60325
$ws.Range("E383").Value = 151.57
$ws.Range("F383").Value = -102
$ws.Range("G383").Value = -12939.72

$ws.Range("B421").Value = 57857
$ws.Range("F421").Value = 3
$ws.Range("G421").Value = 453.51

$ws.Range("B422").Value = 63008
$ws.Range("F422").Value = 504
$ws.Range("G422").Value = 76189.67999999999

$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F431").Value = 36
$ws.Range("G431").Value = 2140.92

$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47

$ws.Range("B579").Value = 53757
$ws.Range("E579").Value = 16.08
$ws.Range("F579").Value = -159
$ws.Range("G579").Value = -2138.55

$ws.Range("B580").Value = 65069
$ws.Range("E580").Value = 14.3
$ws.Range("F580").Value = 172
$ws.Range("G580").Value = 2313.4

$ws.Range("B583").Value = 53263
$ws.Range("E583").Value = 15.29
$ws.Range("F583").Value = -309
$ws.Range("G583").Value = -3958.29

$ws.Range("B584").Value = 65066
$ws.Range("E584").Value = 13.61
$ws.Range("F584").Value = 313
$ws.Range("G584").Value = 4009.53

$ws.Range("B586").Value = 64915
$ws.Range("E586").Value = 20.98
$ws.Range("F586").Value = 40
$ws.Range("G586").Value = 789.2

$ws.Range("B587").Value = 45695
$ws.Range("E587").Value = 23.58
$ws.Range("F587").Value = -36
$ws.Range("G587").Value = -710.28

$ws.Range("B687").Value = 64810
$ws.Range("E687").Value = 291.22
$ws.Range("F687").Value = 7
$ws.Range("G687").Value = 1917.44

$ws.Range("B688").Value = 53319
$ws.Range("E688").Value = 310.64
$ws.Range("F688").Value = -6
$ws.Range("G688").Value = -1643.52

$ws.Range("B709").Value = 60025
$ws.Range("E709").Value = 37.22
$ws.Range("F709").Value = -98
$ws.Range("G709").Value = -3217.34

$ws.Range("B710").Value = 64833
$ws.Range("E710").Value = 34.9
$ws.Range("F710").Value = 99
$ws.Range("G710").Value = 3250.17

$ws.Range("B720").Value = 60022
$ws.Range("E720").Value = 37.22
$ws.Range("F720").Value = -113
$ws.Range("G720").Value = -3709.79

$ws.Range("B721").Value = 64830
$ws.Range("E721").Value = 34.9
$ws.Range("F721").Value = 117
$ws.Range("G721").Value = 3841.11

$ws.Range("B872").Value = 65079
$ws.Range("F872").Value = 21
$ws.Range("G872").Value = 858.27

$ws.Range("B873").Value = 65362
$ws.Range("F873").Value = 2
$ws.Range("G873").Value = 81.73999999999999

